$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (Strike#) values for rows 2-8, replacing the previous
# Strike#-derived figures with the newly calculated K values.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 4
